# Refresh the "cryptos" price/volume table (GitHub Actions scheduled update).
#
# Column D ("Price") holds numeric-looking text such as "94.20" or
# "1.890.42" (a thousands-grouped display string, not a real number).
# Assigning a bare numeric-looking string to Range.Value lets Excel's COM
# layer auto-convert it to a real number, which silently drops trailing
# zeros / switches to scientific notation / reformats thousands separators
# - changing both the stored cell type and, for some values, the displayed
# text. Prefixing the literal with a single quote is exactly what typing
# text into a cell in the Excel UI does: it forces the entry to be stored
# verbatim as text. Resetting the cell's Style to "Normal" afterwards
# clears the transient quote-prefix flag so the cell's style matches the
# source workbook (these data cells carry no explicit style).
#
# Column E ("Volume(1h)") values already contain padding spaces and a "%"
# sign, so Excel leaves them as text automatically - no special handling
# needed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.276.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.83%  '
$ws.Range("D3").Value = "'1.892.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.14%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'323.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.69%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = "'0.5179"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("D8").Value = "'0.4015"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.04%  '
$ws.Range("D9").Value = "'0.08411"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.25%  '
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("D11").Value = "'1.114"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("E12").Value = '  +10.04%  '
$ws.Range("D13").Value = "'6.433"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.06%  '
$ws.Range("D14").Value = "'1.886.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.25%  '
$ws.Range("D15").Value = "'7.316"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("D16").Value = "'1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'94.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = "'0.00001109"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.72%  '
$ws.Range("D19").Value = "'0.06644"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'18.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.32%  '
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").Value = "'5.949"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.43%  '
$ws.Range("D23").Value = "'30.261.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.79%  '
$ws.Range("D24").Value = "'11.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("D25").Value = "'2.231"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.88%  '
$ws.Range("D26").Value = "'2.112.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("E27").Value = '  +2.53%  '
$ws.Range("D28").Value = "'161.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.61%  '
$ws.Range("D29").Value = "'2.334"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.55%  '
$ws.Range("D30").Value = "'129.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.19%  '
$ws.Range("D31").Value = "'1.087"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.33%  '
$ws.Range("E32").Value = '  -0.67%  '
$ws.Range("D33").Value = "'6.092"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.80%  '
$ws.Range("D34").Value = "'3.746"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.43%  '
$ws.Range("D35").Value = "'0.02495"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.30%  '
$ws.Range("D36").Value = "'0.06530"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.73%  '
$ws.Range("D37").Value = "'5.321"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.13%  '
$ws.Range("D38").Value = "'0.2196"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("D39").Value = "'1.220"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.29%  '
$ws.Range("D40").Value = "'8.817"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.70%  '
$ws.Range("D41").Value = "'11.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.42%  '
$ws.Range("D42").Value = "'0.6503"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.76%  '
$ws.Range("D43").Value = "'1.229"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.97%  '
$ws.Range("D44").Value = "'0.6086"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.91%  '
$ws.Range("D45").Value = "'13.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.30%  '
$ws.Range("D46").Value = "'3.682"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.44%  '
$ws.Range("D47").Value = "'2.051"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.90%  '
$ws.Range("D48").Value = "'1.236"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.31%  '
$ws.Range("D49").Value = "'124.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.33%  '
$ws.Range("D50").Value = "'1.160"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.55%  '
$ws.Range("D51").Value = "'79.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.62%  '
